$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Text
$text = $text.Replace("1000 Bs = 2.33 = 8878.33 pesos", "1000 Bs = 2.32 = 8863.11 pesos")
$text = $text.Replace("8878.33 pesos = 2.32 = 955.97 Bs", "8863.11 pesos = 2.31 = 959.36 Bs")
$cellA1.Value = $text

# --- Sheet "tasas": update the rate table values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 431
$wsTasas.Range("O10").Value = 3820
$wsTasas.Range("N12").Value = 3833.99
$wsTasas.Range("O12").Value = 415
